# Applies the "TODOs done, added netlist RTL screenshots" update to the
# FP-log activity-log workbook: fills in the student-info header cells on
# all three sheets, refreshes/extends the Part 2 log descriptions, adds the
# new dated log rows 22-27 (functional/timing sim + RTL netlist screenshots),
# and restores each sheet's last-used selection.

$wb = $excel.ActiveWorkbook

$wsPart1 = $wb.Worksheets.Item("Activity Log - Part 1")
$wsPart2 = $wb.Worksheets.Item("Activity Log - Part 2")
$wsPart3 = $wb.Worksheets.Item("Activity Log - Part 3")

# ---------------------------------------------------------------------------
# Header block (Student Name / Full Student Number / Group Number) is filled
# in on every sheet exactly the same way (Part 3 already had it).
# ---------------------------------------------------------------------------
foreach ($ws in @($wsPart1, $wsPart2)) {
    $ws.Range("B1").Value = "Ruelt Yean (Ryan), Kiew"
    $ws.Range("B2").Value = 301290779
    $ws.Range("B3").Value = "G47"
}

# ---------------------------------------------------------------------------
# "Activity Log - Part 2" is where all the real log data lives.
# ---------------------------------------------------------------------------

# Reworded / expanded descriptions for the already-logged rows 6-21.
$wsPart2.Range("G6").Value  = "Read through Part 2 pdf and set up project folder for initial compilation"
$wsPart2.Range("G7").Value  = "Reviewed Barrel Shifter design implementation - reviewed notes and online sources"
$wsPart2.Range("G8").Value  = "First implementation  of Barrel Shifter - uses 3 different MUX's for simplicity"
$wsPart2.Range("G9").Value  = "Second implementation of Barrel Shifter - attempting to use one MUX entity"
$wsPart2.Range("G10").Value = "Fixing compilation errors for second implementation of Barrel Shifters"
$wsPart2.Range("G11").Value = "Fixing compilation errors for second implementation of Barrel Shifters (DONE)"
$wsPart2.Range("G12").Value = "Third implementaiton of Barrel Shifter - instead of using complicated logic within MUX, just pass options into MUX"
$wsPart2.Range("G13").Value = "Finished implementation of SLL and SRL"
$wsPart2.Range("G14").Value = "Finished implementation of SRA"
$wsPart2.Range("G15").Value = "Reviewed group's code; fixed compilation errors"
$wsPart2.Range("G16").Value = "Fixed logic errors for SRA, SLL and SRL"
$wsPart2.Range("G17").Value = "Fixed logic errors for ShiftUnit"
$wsPart2.Range("G18").Value = "Fixed logic errors for ShiftUnit (DONE) - implemented 32-bit shifting for a shift of more than 32 bits"
$wsPart2.Range("G19").Value = "Fixed some errors that showed up in our modified ArithUnit.vhd - wrong parameters were being passed"
$wsPart2.Range("G20").Value = "Fixed logic errors for ExecUnit"
$wsPart2.Range("G21").Value = "Fixed logic errors in ExecUnit and ShiftUnit, now passes all test cases (DONE)"

# Rows 14-21 already had their date/time entries but were missing the
# last-4-digits (B) and date (C) columns - fill those in now.
$rowInfo = @{
    14 = @{ B = 779; C = 43933 }
    15 = @{ B = 779; C = 43933 }
    16 = @{ B = 779; C = 43934 }
    17 = @{ B = 779; C = 43934 }
    18 = @{ B = 779; C = 43934 }
    19 = @{ B = 779; C = 43935 }
    20 = @{ B = 779; C = 43935 }
    21 = @{ B = 779; C = 43935 }
}
foreach ($r in $rowInfo.Keys) {
    $info = $rowInfo[$r]
    $wsPart2.Range("B$r").Value = $info.B
    $wsPart2.Range("C$r").Value = $info.C
}

# Brand-new rows 22-27: functional/timing simulation screenshots + RTL
# netlist viewer screenshots for ShiftUnit and ExecUnit.
$wsPart2.Range("B22").Value = 779
$wsPart2.Range("C22").Value = 43936
$wsPart2.Range("D22").Value = 0.84375
$wsPart2.Range("E22").Value = 0.88124999999999998
$wsPart2.Range("G22").Value = "Compiling functional simulation screenshots for ShiftUnit"

$wsPart2.Range("B23").Value = 779
$wsPart2.Range("C23").Value = 43936
$wsPart2.Range("D23").Value = 0.89583333333333337
$wsPart2.Range("E23").Value = 0.93055555555555547
$wsPart2.Range("G23").Value = "Compiling timing simulation screenshots for ShiftUnit"

$wsPart2.Range("B24").Value = 779
$wsPart2.Range("C24").Value = 43937
$wsPart2.Range("D24").Value = 0.34930555555555554
$wsPart2.Range("E24").Value = 0.39374999999999999
$wsPart2.Range("G24").Value = "Compiling timing simulation screenshots for ShiftUnit (DONE)"

$wsPart2.Range("B25").Value = 779
$wsPart2.Range("C25").Value = 43937
$wsPart2.Range("D25").Value = 0.39374999999999999
$wsPart2.Range("E25").Value = 0.4055555555555555
$wsPart2.Range("G25").Value = "Compiling functional and timing simulation screenshots for ExecUnit "

$wsPart2.Range("B26").Value = 779
$wsPart2.Range("C26").Value = 43937
$wsPart2.Range("D26").Value = 0.5854166666666667
$wsPart2.Range("E26").Value = 0.59861111111111109
$wsPart2.Range("G26").Value = "Adding comments to code for clarification"

$wsPart2.Range("B27").Value = 779
$wsPart2.Range("C27").Value = 43937
$wsPart2.Range("D27").Value = 0.59861111111111109
$wsPart2.Range("E27").Value = 0.61527777777777781
$wsPart2.Range("G27").Value = "Compiling RTL netlost viewer screenshots for ShiftUnit and ExecUnit"

# ---------------------------------------------------------------------------
# Restore each sheet's last active selection / view position.
# ---------------------------------------------------------------------------
$wsPart1.Activate()
$wsPart1.Range("D6").Select()

$wsPart2.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 16
    $excel.ActiveWindow.ScrollColumn = 1
    $excel.ActiveWindow.Zoom = 100
} catch {
}
$wsPart2.Range("D32").Select()
